$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: 2021 year, fully populated ---
$ws.Range("A7").Value = "2021年"
$ws.Range("B7").Value = 104.1
$ws.Range("C7").Value = 101.5
$ws.Range("D7").Value = 117.1
$ws.Range("E7").Value = 99.5
$ws.Range("F7").Value = 105.3
$ws.Range("G7").Value = 102.3
$ws.Range("H7").Value = 100.9

# --- Row 8: 2022 year, only the first indicator (column B) has been published so far ---
$ws.Range("A8").Value = "2022年"
$ws.Range("B8").Value = 105.2

# The remaining columns (C:H) for 2022 are still blank placeholders (no data yet),
# but the row is present in the sheet, so materialize those cells with the sheet's
# normal/default style instead of leaving them completely absent.
$ws.Range("C8:H8").Style = "Normal"

# Column A uses the same bold/centered "year label" style as the rows above it.
$ws.Range("A6").Copy()
$ws.Range("A7:A8").PasteSpecial(-4122) | Out-Null
